# ThreeConesMethod.xlsx - update solver inputs (angle 1, radius, step size,
# and the starting "z guess" for the second tracker column) and move the
# selection/viewport over to the I-column guess table, as the author did
# while validating the new low/high-guess support.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# theta1 (B1): 56 -> 45 degrees (C1 = RADIANS(90-B1) recalculates automatically)
$ws.Range("B1").Value = 45

# s, the target separation between trackers (F1): 350 -> 101
$ws.Range("F1").Value = 101

# step size used to build the I-column guesses (I5): 3 -> 0.1
$ws.Range("I5").Value = 0.1

# starting guess for the I-column table (I7): 675 -> 141
$ws.Range("I7").Value = 141

# Move the selection/viewport to the guess table the author was inspecting
$ws.Range("I6").Select()
